$d = $word.ActiveDocument

# Convert each `m:userdoc`/`m:enduserdoc` field (fldChar begin/instrText/fldChar end)
# into a single plain-text run holding the literal "{ ... }" token text, e.g.
#   { fldChar begin }{ instrText " m:userdoc 'zone1' " }{ fldChar end }
# becomes
#   <w:t xml:space="preserve">{m:userdoc 'zone1'}</w:t>
# This mirrors the TokenIteratorFieldRewriterSplit parser switching from
# real Word fields to plain-text "{...}" tokens in the template.

function ConvertTo-LiteralToken([string]$code) {
    # field code text looks like " m:userdoc 'zone1' " or " m:enduserdoc "
    $trimmed = $code.Trim()
    return "{" + $trimmed + "}"
}

function Set-PlainTextRun($range, [string]$text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xmlPkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">' +
        $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xmlPkg)
}

# All fields in this template are the userdoc/enduserdoc markers that need
# rewriting, so repeatedly take the first field, remember where it starts,
# delete it (collapsing it to an empty spot) and drop the literal text run
# in its place.
while ($d.Fields.Count -gt 0) {
    $f = $d.Fields.Item(1)
    $token = ConvertTo-LiteralToken($f.Code.Text)
    $start = $f.Code.Start - 1
    $f.Delete()
    $r = $d.Range($start, $start)
    Set-PlainTextRun $r $token
}
